$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.662.62"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "1.869.63"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").Value = "'331.86"
$ws.Range("E5").Value = "  +3.05%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = "  +4.06%  "

$ws.Range("E8").Value = "  +2.35%  "

$ws.Range("D9").Value = "'47.86"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").Value = "'0.08059"
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("D11").Value = "'1.022"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "'21.77"
$ws.Range("E12").Value = "  +2.17%  "

$ws.Range("D13").Value = "1.873.54"
$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("D14").Value = "'5.936"
$ws.Range("E14").Value = "  +1.14%  "

$ws.Range("D15").Value = "'7.140"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'0.00001047"
$ws.Range("E17").Value = "  +1.66%  "

$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").Value = "'0.06638"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").Value = "'17.13"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").Value = "27.675.27"
$ws.Range("E22").Value = "  +1.35%  "

$ws.Range("D23").Value = "'5.488"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").Value = "'11.00"

$ws.Range("D25").Value = "'2.313"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("D26").Value = "2.089.59"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").Value = "'158.68"
$ws.Range("E27").Value = "  +4.59%  "

$ws.Range("D28").Value = "'20.19"
$ws.Range("E28").Value = "  +2.60%  "

$ws.Range("D29").Value = "'2.090"
$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("D30").Value = "'5.551"
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("E31").Value = "  +1.75%  "

$ws.Range("D32").Value = "'0.9653"
$ws.Range("E32").Value = "  +3.59%  "

$ws.Range("E33").Value = "  +2.26%  "

$ws.Range("D34").Value = "'1.444"
$ws.Range("E34").Value = "  -2.41%  "

$ws.Range("D35").Value = "'3.594"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Value = "'5.313"
$ws.Range("E36").Value = "  +0.91%  "

$ws.Range("E37").Value = "  +1.60%  "

$ws.Range("D38").Value = "'0.06084"
$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("D39").Value = "'1.231"
$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("D40").Value = "'8.124"
$ws.Range("E40").Value = "  -1.72%  "

$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").Value = "'0.5981"
$ws.Range("E42").Value = "  +1.29%  "

$ws.Range("D43").Value = "'0.1894"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("D44").Value = "'10.22"
$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("D45").Value = "'1.252"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "'0.5708"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").Value = "'12.17"
$ws.Range("E47").Value = "  +2.34%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.389"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.936"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").Value = "'0.06852"
$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("D51").Value = "'114.22"
$ws.Range("E51").Value = "  +5.72%  "
